$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values for re-run SGNN dialog act annotation
$ws.Range("I5").Value = "sv"
$ws.Range("J5").Value = "Statement-opinion"
$ws.Range("I6").Value = "b"
$ws.Range("J6").Value = "Acknowledge (Backchannel)"
$ws.Range("I7").Value = "sv"
$ws.Range("J7").Value = "Statement-opinion"
$ws.Range("I11").Value = "b"
$ws.Range("J11").Value = "Acknowledge (Backchannel)"
$ws.Range("I18").Value = "aa"
$ws.Range("J18").Value = "Agree/Accept"
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I26").Value = "b"
$ws.Range("J26").Value = "Acknowledge (Backchannel)"
$ws.Range("I30").Value = "sd"
$ws.Range("J30").Value = "Statement-non-opinion"
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"
$ws.Range("I45").Value = "%"
$ws.Range("J45").Value = "Uninterpretable"
$ws.Range("I47").Value = "b"
$ws.Range("J47").Value = "Acknowledge (Backchannel)"
$ws.Range("I54").Value = "aa"
$ws.Range("J54").Value = "Agree/Accept"
$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"
$ws.Range("I80").Value = "sd"
$ws.Range("J80").Value = "Statement-non-opinion"
$ws.Range("I81").Value = "sd"
$ws.Range("J81").Value = "Statement-non-opinion"
$ws.Range("I103").Value = "ba"
$ws.Range("J103").Value = "Appreciation"
$ws.Range("I106").Value = "sd"
$ws.Range("J106").Value = "Statement-non-opinion"
$ws.Range("I118").Value = "b"
$ws.Range("J118").Value = "Acknowledge (Backchannel)"
$ws.Range("I122").Value = "%"
$ws.Range("J122").Value = "Uninterpretable"
$ws.Range("I132").Value = "b"
$ws.Range("J132").Value = "Acknowledge (Backchannel)"
$ws.Range("I139").Value = "aa"
$ws.Range("J139").Value = "Agree/Accept"
$ws.Range("I141").Value = "sv"
$ws.Range("J141").Value = "Statement-opinion"
$ws.Range("I143").Value = "sd"
$ws.Range("J143").Value = "Statement-non-opinion"
$ws.Range("I162").Value = "aa"
$ws.Range("J162").Value = "Agree/Accept"
$ws.Range("I196").Value = "aa"
$ws.Range("J196").Value = "Agree/Accept"
$ws.Range("I210").Value = "aa"
$ws.Range("J210").Value = "Agree/Accept"
$ws.Range("I211").Value = "aa"
$ws.Range("J211").Value = "Agree/Accept"
$ws.Range("I224").Value = "sd"
$ws.Range("J224").Value = "Statement-non-opinion"
$ws.Range("I245").Value = "sd"
$ws.Range("J245").Value = "Statement-non-opinion"
$ws.Range("I256").Value = "%"
$ws.Range("J256").Value = "Uninterpretable"
$ws.Range("I260").Value = "b"
$ws.Range("J260").Value = "Acknowledge (Backchannel)"
$ws.Range("I264").Value = "b"
$ws.Range("J264").Value = "Acknowledge (Backchannel)"
$ws.Range("I269").Value = "sv"
$ws.Range("J269").Value = "Statement-opinion"
$ws.Range("I271").Value = "b"
$ws.Range("J271").Value = "Acknowledge (Backchannel)"
$ws.Range("I273").Value = "b"
$ws.Range("J273").Value = "Acknowledge (Backchannel)"
$ws.Range("I276").Value = "b"
$ws.Range("J276").Value = "Acknowledge (Backchannel)"
$ws.Range("I278").Value = "aa"
$ws.Range("J278").Value = "Agree/Accept"
$ws.Range("I293").Value = "aa"
$ws.Range("J293").Value = "Agree/Accept"
